$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on D-column cells whose new price strings would otherwise
# be auto-parsed by Excel as numbers (single decimal point), to preserve them
# as text values, matching the original inlineStr cell type.
foreach ($r in @(5,8,10,16,17,19,23,24,25,26,27,34,35,37,39,44,47)) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.995.47"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "1.678.34"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "215.30"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("E6").Value = "  +1.26%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "0.251"
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("E9").Value = "  +0.40%  "
$ws.Range("D10").Value = "20.38"
$ws.Range("E10").Value = "  +1.48%  "
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("D12").Value = "1.915.43"
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("D13").Value = "1.680.37"
$ws.Range("E13").Value = "  +0.50%  "
$ws.Range("E14").Value = "  +0.22%  "
$ws.Range("E15").Value = "  +1.70%  "
$ws.Range("D16").Value = "65.72"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("D17").Value = "8.24"
$ws.Range("E17").Value = "  +6.79%  "
$ws.Range("D18").Value = "27.025.34"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").Value = "235.95"
$ws.Range("E19").Value = "  +0.42%  "
$ws.Range("D20").Value = "0.0₃0734"
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D23").Value = "9.19"
$ws.Range("E23").Value = "  -0.48%  "
$ws.Range("D24").Value = "2.17"
$ws.Range("E24").Value = "  -2.90%  "
$ws.Range("D25").Value = "146.01"
$ws.Range("E25").Value = "  +0.47%  "
$ws.Range("D26").Value = "7.23"
$ws.Range("E26").Value = "  +1.09%  "
$ws.Range("D27").Value = "16.07"
$ws.Range("E27").Value = "  +1.25%  "
$ws.Range("E28").Value = "  -1.57%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("E31").Value = "  -0.52%  "
$ws.Range("E32").Value = "  -0.08%  "
$ws.Range("D33").Value = "1.480.49"
$ws.Range("E33").Value = "  +1.86%  "
$ws.Range("D34").Value = "3.18"
$ws.Range("E34").Value = "  +1.21%  "
$ws.Range("D35").Value = "1.68"
$ws.Range("E35").Value = "  +5.03%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").Value = "0.583"
$ws.Range("E37").Value = "  +2.60%  "
$ws.Range("E38").Value = "  +2.60%  "
$ws.Range("D39").Value = "0.904"
$ws.Range("E39").Value = "  +1.31%  "
$ws.Range("E40").Value = "  -3.89%  "
$ws.Range("E41").Value = "  +0.95%  "
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("E43").Value = "  +1.98%  "
$ws.Range("D44").Value = "67.45"
$ws.Range("E44").Value = "  +2.59%  "
$ws.Range("D45").Value = "1.819.47"
$ws.Range("E45").Value = "  +0.10%  "
$ws.Range("E46").Value = "  +0.43%  "
$ws.Range("D47").Value = "90.46"
$ws.Range("E47").Value = "  +0.24%  "
$ws.Range("E48").Value = "  +0.76%  "
$ws.Range("E49").Value = "  -0.28%  "
$ws.Range("E50").Value = "  +1.50%  "
$ws.Range("E51").Value = "  +0.09%  "
